$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Try setting explicit protection via .Locked (already true since xf2 has locked=1, hidden=0)
$ws.Range("D4").Locked = $true
$ws.Range("D4").FormulaHidden = $false
$ws.Range("D4").Value = "motor_id"
